$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: likes count update
$ws.Range("H2").Value = 312

# Row 5
$ws.Range('B5').Value = 'Snowflake'
$ws.Range('C5').Value = '{''avatarUrl'': ''https://cdn-avatars.huggingface.co/v1/production/uploads/64ba2f59a6ccf0f64b4ad254/eTDA37yFwUVP45c1WTSs2.png'', ''fullname'': ''Snowflake'', ''name'': ''Snowflake'', ''type'': ''org'', ''isHf'': False, ''isEnterprise'': False}'
$ws.Range('D5').Value = 33257
$ws.Range('F5').Value = 'Snowflake/snowflake-arctic-embed-l'
$ws.Range('G5').Value = '2024-04-18T19:58:11.000Z'
$ws.Range('H5').Value = 59
$ws.Range('I5').Value = 'sentence-similarity'
$ws.Range('M5').Value = '1.34GB | 299MB | 669MB | 337MB | 318MB | 337MB | 337MB'
$ws.Range('N5').Value = 313524224

# Row 6
$ws.Range('B6').Value = 'WhereIsAI'
$ws.Range('C6').Value = '{''avatarUrl'': ''https://www.gravatar.com/avatar/e81bd32cb5ee88835824ad6b60d05697?d=retro&size=100'', ''fullname'': ''WhereIsAI'', ''name'': ''WhereIsAI'', ''type'': ''org'', ''isHf'': False, ''isEnterprise'': False}'
$ws.Range('D6').Value = 277613
$ws.Range('F6').Value = 'WhereIsAI/UAE-Large-V1'
$ws.Range('G6').Value = '2024-05-03T02:31:54.000Z'
$ws.Range('H6').Value = 177
$ws.Range('I6').Value = 'feature-extraction'
$ws.Range('M6').Value = '1.34GB | 669MB | 337MB'
$ws.Range('N6').Value = 353370112

# Row 7
$ws.Range('B7').Value = 'Alibaba-NLP'
$ws.Range('C7').Value = '{''avatarUrl'': ''https://www.gravatar.com/avatar/1ae3fd9f5b9356f196c997d93eb23038?d=retro&size=100'', ''fullname'': ''Alibaba-NLP'', ''name'': ''Alibaba-NLP'', ''type'': ''org'', ''isHf'': False, ''isEnterprise'': False}'
$ws.Range('D7').Value = 75109
$ws.Range('F7').Value = 'Alibaba-NLP/gte-large-en-v1.5'
$ws.Range('G7').Value = '2024-04-26T13:51:26.000Z'
$ws.Range('H7').Value = 57
$ws.Range('I7').Value = 'sentence-similarity'
$ws.Range('M7').Value = '1.75GB | 361MB | 873MB | 446MB | 387MB | 446MB | 446MB'
$ws.Range('N7').Value = 378535936
$ws.Range('O7').Value = 1879048192

# Row 8
$ws.Range('B8').Value = 'Snowflake'
$ws.Range('C8').Value = '{''avatarUrl'': ''https://cdn-avatars.huggingface.co/v1/production/uploads/64ba2f59a6ccf0f64b4ad254/eTDA37yFwUVP45c1WTSs2.png'', ''fullname'': ''Snowflake'', ''name'': ''Snowflake'', ''type'': ''org'', ''isHf'': False, ''isEnterprise'': False}'
$ws.Range('D8').Value = 19315
$ws.Range('F8').Value = 'Snowflake/snowflake-arctic-embed-m'
$ws.Range('G8').Value = '2024-04-18T19:50:37.000Z'
$ws.Range('H8').Value = 63
$ws.Range('I8').Value = 'sentence-similarity'
$ws.Range('M8').Value = '436MB | 144MB | 218MB | 110MB | 149MB | 110MB | 110MB'
$ws.Range('N8').Value = 115343360
$ws.Range('O8').Value = 457179136
